$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 50, shifting existing rows 50-75 down to 51-76.
$ws.Rows.Item(50).Insert()

# Copy the fixed columns from row 51 (the row that was row 50 before the insert)
# into the new row 50, then set the new record's specific values.
$ws.Range("A50").Value = $ws.Range("A51").Value2
$ws.Range("B50").Value = $ws.Range("B51").Value2
$ws.Range("C50").Value = $ws.Range("C51").Value2
$ws.Range("D50").NumberFormat = $ws.Range("D51").NumberFormat
$ws.Range("D50").Value = 44572
$ws.Range("E50").Value = $ws.Range("E51").Value2
$ws.Range("F50").Value = $ws.Range("F51").Value2
$ws.Range("G50").Value = $ws.Range("G51").Value2
$ws.Range("H50").Value = $ws.Range("H51").Value2
$ws.Range("I50").Value = $ws.Range("I51").Value2
$ws.Range("J50").Value = 80
$ws.Range("K50").Value = 23000
$ws.Range("L50").Value = 23000
$ws.Range("M50").Value = 23000
$ws.Range("N50").Value = $ws.Range("N51").Value2
$ws.Range("O50").Value = "Región de La Araucanía"
$ws.Range("P50").Value = 920
$ws.Range("Q50").Value = $ws.Range("Q51").Value2
$ws.Range("R50").Value = $ws.Range("R51").Value2
